# "Cubiertas horas de octubre" — fill in the October hours-worked row (row 11,
# "Hours worked on this project") on the TimeSheet. Columns B..AF hold days
# 1..31; only the weekday columns carry values (weekend columns stay blank).
# The dependent SUM formulas (row 18 daily totals, AE23/AE30 monthly totals)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$hours = @{
    "B11" = 9.5
    "C11" = 10
    "D11" = 9
    "E11" = 2.5
    "H11" = 8.25
    "I11" = 9.5
    "J11" = 9
    "K11" = 9
    "L11" = 3
    "O11" = 9.5
    "P11" = 8.5
    "Q11" = 10.5
    "S11" = 4.5
    "V11" = 8.5
    "W11" = 8.5
    "X11" = 5.5
    "Z11" = 7.5
    "AC11" = 10
    "AD11" = 9.5
}

foreach ($addr in $hours.Keys) {
    $ws.Range($addr).Value = $hours[$addr]
}

# Last cell touched while entering the October data.
$ws.Range("AI24").Select() | Out-Null
